$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 195.435389
$ws.Range("H2").Value = 586.306167
$ws.Range("I2").Value = 0.3095741734129938
$ws.Range("J2").Value = 0.3095741734129938
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.98286266666667
$ws.Range("N2").Value = 44.948588
$ws.Range("O2").Value = 0.1958132590302862
$ws.Range("P2").Value = 0.1958132590302861
$ws.Range("Q2").Value = 2928.181593593577
$ws.Range("R2").Value = 26353.63434234219
$ws.Range("S2").Value = 0.06061872780760529
$ws.Range("T2").Value = 0.06061872780760526

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 195.435389
$ws.Range("H3").Value = 586.306167
$ws.Range("I3").Value = 0.3095741734129938
$ws.Range("J3").Value = 0.3095741734129938
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 53.32328933333334
$ws.Range("N3").Value = 159.969868
$ws.Range("O3").Value = 0.6968899935126925
$ws.Range("P3").Value = 0.6968899935126924
$ws.Range("Q3").Value = 10421.25779361955
$ws.Range("R3").Value = 93791.32014257596
$ws.Range("S3").Value = 0.2157391437014784
$ws.Range("T3").Value = 0.2157391437014783

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 195.435389
$ws.Range("H4").Value = 586.306167
$ws.Range("I4").Value = 0.3095741734129938
$ws.Range("J4").Value = 0.3095741734129938
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.209926333333334
$ws.Range("N4").Value = 24.629779
$ws.Range("O4").Value = 0.1072967474570214
$ws.Range("P4").Value = 0.1072967474570214
$ws.Range("Q4").Value = 1604.510146616344
$ws.Range("R4").Value = 14440.59131954709
$ws.Range("S4").Value = 0.03321630190391014
$ws.Range("T4").Value = 0.03321630190391014

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 79.82725266666667
$ws.Range("H5").Value = 239.481758
$ws.Range("I5").Value = 0.1264482133280045
$ws.Range("J5").Value = 0.1264482133280045
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.98286266666667
$ws.Range("N5").Value = 44.948588
$ws.Range("O5").Value = 0.1958132590302862
$ws.Range("P5").Value = 0.1958132590302861
$ws.Range("Q5").Value = 1196.040763761967
$ws.Range("R5").Value = 10764.3668738577
$ws.Range("S5").Value = 0.02476023675031342
$ws.Range("T5").Value = 0.02476023675031342

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 79.82725266666667
$ws.Range("H6").Value = 239.481758
$ws.Range("I6").Value = 0.1264482133280045
$ws.Range("J6").Value = 0.1264482133280045
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 53.32328933333334
$ws.Range("N6").Value = 159.969868
$ws.Range("O6").Value = 0.6968899935126925
$ws.Range("P6").Value = 0.6968899935126924
$ws.Range("Q6").Value = 4256.651690629772
$ws.Range("R6").Value = 38309.86521566795
$ws.Range("S6").Value = 0.08812049456584461
$ws.Range("T6").Value = 0.08812049456584457

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 79.82725266666667
$ws.Range("H7").Value = 239.481758
$ws.Range("I7").Value = 0.1264482133280045
$ws.Range("J7").Value = 0.1264482133280045
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.209926333333334
$ws.Range("N7").Value = 24.629779
$ws.Range("O7").Value = 0.1072967474570214
$ws.Range("P7").Value = 0.1072967474570214
$ws.Range("Q7").Value = 655.3758637857203
$ws.Range("R7").Value = 5898.382774071482
$ws.Range("S7").Value = 0.01356748201184646
$ws.Range("T7").Value = 0.01356748201184646

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 162.5116576666667
$ws.Range("H8").Value = 487.534973
$ws.Range("I8").Value = 0.2574222219914007
$ws.Range("J8").Value = 0.2574222219914007
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.98286266666667
$ws.Range("N8").Value = 44.948588
$ws.Range("O8").Value = 0.1958132590302862
$ws.Range("P8").Value = 0.1958132590302861
$ws.Range("Q8").Value = 2434.889848552014
$ws.Range("R8").Value = 21914.00863696812
$ws.Range("S8").Value = 0.05040668423495397
$ws.Range("T8").Value = 0.05040668423495396

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 162.5116576666667
$ws.Range("H9").Value = 487.534973
$ws.Range("I9").Value = 0.2574222219914007
$ws.Range("J9").Value = 0.2574222219914007
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 53.32328933333334
$ws.Range("N9").Value = 159.969868
$ws.Range("O9").Value = 0.6968899935126925
$ws.Range("P9").Value = 0.6968899935126924
$ws.Range("Q9").Value = 8665.656141799285
$ws.Range("R9").Value = 77990.90527619357
$ws.Range("S9").Value = 0.1793949706136101
$ws.Range("T9").Value = 0.1793949706136101

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 162.5116576666667
$ws.Range("H10").Value = 487.534973
$ws.Range("I10").Value = 0.2574222219914007
$ws.Range("J10").Value = 0.2574222219914007
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.209926333333334
$ws.Range("N10").Value = 24.629779
$ws.Range("O10").Value = 0.1072967474570214
$ws.Range("P10").Value = 0.1072967474570214
$ws.Range("Q10").Value = 1334.208737751218
$ws.Range("R10").Value = 12007.87863976097
$ws.Range("S10").Value = 0.02762056714283662
$ws.Range("T10").Value = 0.02762056714283662

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 159.7910413333334
$ws.Range("H11").Value = 479.3731240000001
$ws.Range("I11").Value = 0.2531127028358626
$ws.Range("J11").Value = 0.2531127028358626
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.98286266666667
$ws.Range("N11").Value = 44.948588
$ws.Range("O11").Value = 0.1958132590302862
$ws.Range("P11").Value = 0.1958132590302861
$ws.Range("Q11").Value = 2394.127227660991
$ws.Range("R11").Value = 21547.14504894891
$ws.Range("S11").Value = 0.04956282324425462
$ws.Range("T11").Value = 0.04956282324425459

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 159.7910413333334
$ws.Range("H12").Value = 479.3731240000001
$ws.Range("I12").Value = 0.2531127028358626
$ws.Range("J12").Value = 0.2531127028358626
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 53.32328933333334
$ws.Range("N12").Value = 159.969868
$ws.Range("O12").Value = 0.6968899935126925
$ws.Range("P12").Value = 0.6968899935126924
$ws.Range("Q12").Value = 8520.583929891962
$ws.Range("R12").Value = 76685.25536902765
$ws.Range("S12").Value = 0.1763917098372644
$ws.Range("T12").Value = 0.1763917098372643

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 159.7910413333334
$ws.Range("H13").Value = 479.3731240000001
$ws.Range("I13").Value = 0.2531127028358626
$ws.Range("J13").Value = 0.2531127028358626
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.209926333333334
$ws.Range("N13").Value = 24.629779
$ws.Range("O13").Value = 0.1072967474570214
$ws.Range("P13").Value = 0.1072967474570214
$ws.Range("Q13").Value = 1311.872678073289
$ws.Range("R13").Value = 11806.8541026596
$ws.Range("S13").Value = 0.02715816975434365
$ws.Range("T13").Value = 0.02715816975434365

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 33.73857866666666
$ws.Range("H14").Value = 101.215736
$ws.Range("I14").Value = 0.05344268843173843
$ws.Range("J14").Value = 0.05344268843173842
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 14.98286266666667
$ws.Range("N14").Value = 44.948588
$ws.Range("O14").Value = 0.1958132590302862
$ws.Range("P14").Value = 0.1958132590302861
$ws.Range("Q14").Value = 505.5004907311964
$ws.Range("R14").Value = 4549.504416580768
$ws.Range("S14").Value = 0.01046478699315887
$ws.Range("T14").Value = 0.01046478699315887

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 33.73857866666666
$ws.Range("H15").Value = 101.215736
$ws.Range("I15").Value = 0.05344268843173843
$ws.Range("J15").Value = 0.05344268843173842
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 53.32328933333334
$ws.Range("N15").Value = 159.969868
$ws.Range("O15").Value = 0.6968899935126925
$ws.Range("P15").Value = 0.6968899935126924
$ws.Range("Q15").Value = 1799.051991938094
$ws.Range("R15").Value = 16191.46792744285
$ws.Range("S15").Value = 0.03724367479449504
$ws.Range("T15").Value = 0.03724367479449504

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 33.73857866666666
$ws.Range("H16").Value = 101.215736
$ws.Range("I16").Value = 0.05344268843173843
$ws.Range("J16").Value = 0.05344268843173842
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.209926333333334
$ws.Range("N16").Value = 24.629779
$ws.Range("O16").Value = 0.1072967474570214
$ws.Range("P16").Value = 0.1072967474570214
$ws.Range("Q16").Value = 276.9912454447049
$ws.Range("R16").Value = 2492.921209002344
$ws.Range("S16").Value = 0.005734226644084517
$ws.Range("T16").Value = 0.005734226644084516
